# renaming @lvl.core to @lvl/front-end, in accordance with
# https://github.com/jonathanburrows/lavalav/issues/7
#
# Backlog housekeeping on the "backlog" sheet:
#   - POC - Windows Federated STS        -> Deferred
#   - POC - Developer management         -> Deferred
#   - .gitignore changes                 -> Completed
#   - Bitwise Architecture Consistency   -> Completed
#   - .net Framework Consistency         -> Completed
#   - Change folder for @lvl.core to front-end -> In Progress (this is the
#     task the commit itself is doing)
#   - Fixing dotnet runner               -> Completed, with a description
#     explaining the custom dotnet runner was removed in favor of the
#     built-in one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("backlog")

$ws.Range("E33").Value = "Deferred"
$ws.Range("E34").Value = "Deferred"

$ws.Range("E35").Value = "Completed"
$ws.Range("E36").Value = "Completed"
$ws.Range("E37").Value = "Completed"

$ws.Range("E38").Value = "In Progress"

$ws.Range("E51").Value = "Completed"
$ws.Range("G51").Value = "Remove the custom dotnet runner, and replace it with the built in one, so that mantenance is decreased."

# Move the view / selection to where the work happened.
[void]$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 35
$win.ScrollColumn = 1
[void]$ws.Range("D50").Select()
